$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Increment the "Inscritos" (column E) count by 1 for the specified rows,
# reflecting one additional registration recorded in each of these rows.
$rows = @(15, 38, 56, 64, 65, 69, 76, 89)

foreach ($r in $rows) {
    $cell = $ws.Range("E$r")
    $cell.Value = $cell.Value2 + 1
}
